$d = $word.ActiveDocument

# Locate the paragraph run that currently reads:
#   "Goodyear, 85338 | Willing to relocate "
# and split it into three runs (same rFonts formatting on each):
#   "Goodyear, "  +  "AZ "  +  "85338 | Willing to relocate "
$found = $d.Content
$found.Find.Execute("Goodyear, 85338 | Willing to relocate ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

# Re-wrap the located span in a fresh Range before editing it; InsertXML
# needs a plain Range (not one still tied to a live Find operation).
$target = $d.Range($found.Start, $found.End)

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
         '<w:r>' + $rPr + '<w:t xml:space="preserve">Goodyear, </w:t></w:r>' + `
         '<w:r>' + $rPr + '<w:t xml:space="preserve">AZ </w:t></w:r>' + `
         '<w:r>' + $rPr + '<w:t xml:space="preserve">85338 | Willing to relocate </w:t></w:r>' + `
       '</w:p>'

$target.InsertXML($xml)
